$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Retailers Format")

# --- Row 2: DSR-0350 / Lalpur retailer -> now "Sinja Telecom" / Md Sanowar Hosen(Sujon)
$ws.Range("C2").Value = "Sinja Telecom"
$ws.Range("E2").Value = "Md Sanowar Hosen(Sujon)"
$ws.Range("I2").Value = "Md Sanowar Hosen(Sujon)"
$ws.Range("J2").Value = 1723771230
$ws.Range("P2").Value = 1723771230
$ws.Range("T2").Value = 1723771230

# --- Row 3: was DSR-0248 / Gurudaspur -> now DSR-0349 / Bagha, Rajshahi "Anika Telecom"
$ws.Range("B3").Value = "DSR-0349"
$ws.Range("C3").Value = "Anika Telecom"
$ws.Range("D3").Value = "Bagha"
$ws.Range("E3").Value = "Md Anisur Rahman(Azad)"
$ws.Range("I3").Value = "Md Anisur Rahman(Azad)"
$ws.Range("J3").Value = 1748937325
$ws.Range("K3").Value = "Rajshahi"
$ws.Range("L3").Value = "Bagha"
$ws.Range("N3").Value = "Monigram Bazar, Bagha, Rajshahi."
$ws.Range("P3").Value = 1748937325
$ws.Range("T3").Value = 1748937325

# --- Row 4: was DSR-0248 / Baraigram -> now DSR-0619 / Lalpur "Fatema Telecom 2"
$ws.Range("B4").Value = "DSR-0619"
$ws.Range("C4").Value = "Fatema Telecom 2"
$ws.Range("D4").Value = "Lalpur"
$ws.Range("E4").Value = "Md Kutub Uddin"
$ws.Range("I4").Value = "Md Kutub Uddin"
$ws.Range("J4").Value = 1719132820
$ws.Range("L4").Value = "Lalpur"
$ws.Range("N4").Value = "Oalia Bazar, Lalpur, Natore."
$ws.Range("P4").Value = 1719132820
$ws.Range("T4").Value = 1719132820

# --- Update the saved cursor/selection state on the sheet view.
$ws.Range("I17:I18").Select()
